$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. EWAR sheet: insert 4 new risk-assessment questions before the existing
#    "organisation responding" question (old row 43), i.e. at row 43.
#    The rows below (old 43-54) are pushed down to 47-58 automatically.
# ---------------------------------------------------------------------------
$ewar = $wb.Worksheets.Item("EWAR")

$ewar.Rows.Item(43).Insert()
$ewar.Rows.Item(43).Insert()
$ewar.Rows.Item(43).Insert()
$ewar.Rows.Item(43).Insert()

# Row 43 - maximum_size
$ewar.Range("A43").Value = "select_one hml"
$ewar.Range("B43").Value = "maximum_size"
$ewar.Range("C43").Value = "Likely maximum size of an outbreak linked with this event"
$ewar.Range("D43").Value = "Taille maximale probable d'une épidémie liée à cet événement"
$ewar.Rows.Item(43).RowHeight = 28.5

# Row 44 - maximum_impact
$ewar.Range("A44").Value = "select_one hml"
$ewar.Range("B44").Value = "maximum_impact"
$ewar.Range("C44").Value = "Likely maximum impact of an outbreak linked with this event"
$ewar.Range("D44").Value = "Impact maximal probable d'une épidémie liée à cet événement"
$ewar.Rows.Item(44).RowHeight = 28.5

# Row 45 - duration_outbreak
$ewar.Range("A45").Value = "select_one time"
$ewar.Range("B45").Value = "duration_outbreak"
$ewar.Range("C45").Value = "Likely duration of an outbreak linked with this event"
$ewar.Range("D45").Value = "Durée probable d'une épidémie liée à cet événement"
$ewar.Rows.Item(45).RowHeight = 28.5

# Row 46 - eprep
$ewar.Range("A46").Value = "select_one ynu"
$ewar.Range("B46").Value = "eprep"
$ewar.Range("C46").Value = "Event included among EPREP scenarios"
$ewar.Range("D46").Value = "Événement inclus dans les scénarios EPREP"

# ---------------------------------------------------------------------------
# 2. EWAR_options sheet: append the choice lists for the new "hml" and
#    "time" select_one questions (rows 55-60).
# ---------------------------------------------------------------------------
$ewarOptions = $wb.Worksheets.Item("EWAR_options")

$ewarOptions.Range("A55").Value = "hml"
$ewarOptions.Range("B55").Value = "high"
$ewarOptions.Range("C55").Value = "High"
$ewarOptions.Range("D55").Value = "Haut"

$ewarOptions.Range("A56").Value = "hml"
$ewarOptions.Range("B56").Value = "medium"
$ewarOptions.Range("C56").Value = "Medium"
$ewarOptions.Range("D56").Value = "Moyen"

$ewarOptions.Range("A57").Value = "hml"
$ewarOptions.Range("B57").Value = "low"
$ewarOptions.Range("C57").Value = "Low"
$ewarOptions.Range("D57").Value = "Bas"

$ewarOptions.Range("A58").Value = "time"
$ewarOptions.Range("B58").Value = "weeks"
$ewarOptions.Range("C58").Value = "Weeks"
$ewarOptions.Range("D58").Value = "Semaines"

$ewarOptions.Range("A59").Value = "time"
$ewarOptions.Range("B59").Value = "months"
$ewarOptions.Range("C59").Value = "Months"
$ewarOptions.Range("D59").Value = "Mois"

$ewarOptions.Range("A60").Value = "time"
$ewarOptions.Range("B60").Value = "years"
$ewarOptions.Range("C60").Value = "Years"
$ewarOptions.Range("D60").Value = "Années"

# ---------------------------------------------------------------------------
# 3. Vaccination_long sheet: row 67 loses its custom row height (25.5 -> the
#    sheet's default of 12.75).
# ---------------------------------------------------------------------------
$vaccLong = $wb.Worksheets.Item("Vaccination_long")
$vaccLong.Rows.Item(67).RowHeight = 12.75

# ---------------------------------------------------------------------------
# 4. Restore view state: scroll position / selection on EWAR and
#    EWAR_options, and make EWAR_options the active tab (as in the edited
#    workbook).
# ---------------------------------------------------------------------------
$ewar.Range("A51").Select()
$ewar.Application.ActiveWindow.ScrollRow = 34

$ewarOptions.Select()
$ewarOptions.Range("F55").Select()
$ewarOptions.Application.ActiveWindow.ScrollRow = 34
